# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45208
$ws.Range("C3").Value = 45208
$ws.Range("C4").Value = 45208
$ws.Range("C5").Value = 45208
